# Atualizado por script em 11-11-2023 20:45
#
# The source scrape re-pulled this round of matches and the row order for
# three pairs of fixtures came back swapped relative to the previous
# export. Also a newly played fixture (Altay x Eyupspor) was appended at
# the bottom of the sheet. Swap the F:V ("home" .. "url_partida") payload
# for each pair while leaving A:E (Indice/pais/torneio/temporada/
# data_partida) untouched, then append the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchRows($rowA, $rowB) {
    $dataA = $ws.Range("F$rowA`:V$rowA").Value2
    $dataB = $ws.Range("F$rowB`:V$rowB").Value2
    $ws.Range("F$rowA`:V$rowA").Value2 = $dataB
    $ws.Range("F$rowB`:V$rowB").Value2 = $dataA
}

# Rows 74 / 75 had their match data swapped back.
Swap-MatchRows 74 75

# Rows 84 / 85 had their match data swapped back.
Swap-MatchRows 84 85

# Rows 89 / 90 had their match data swapped back.
Swap-MatchRows 89 90

# Append new row 105 (Indice 104): Altay 1-7 Eyupspor.
$newRow = 105

# Copy number formats/styles from the row above for the styled columns
# (A -> bold/bordered index style, E -> datetime style) before writing
# the values, so no new duplicate style gets minted.
$ws.Cells.Item($newRow - 1, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)
$ws.Cells.Item($newRow - 1, 5).Copy()
$ws.Cells.Item($newRow, 5).PasteSpecial(-4122)

$ws.Cells.Item($newRow, 1).Value = 104
$ws.Cells.Item($newRow, 2).Value = "turkey"
$ws.Cells.Item($newRow, 3).Value = "1-lig"
$ws.Cells.Item($newRow, 4).Value = "2023-2024"
$ws.Cells.Item($newRow, 5).Value2 = 45241.70833333334
$ws.Cells.Item($newRow, 6).Value = "Altay"
$ws.Cells.Item($newRow, 7).Value = 1
$ws.Cells.Item($newRow, 8).Value = "Eyupspor"
$ws.Cells.Item($newRow, 9).Value = 7
$ws.Cells.Item($newRow, 10).Value = 6.63
$ws.Cells.Item($newRow, 11).Value = "05/11/2023 17:12"
$ws.Cells.Item($newRow, 12).Value = 9.33
$ws.Cells.Item($newRow, 13).Value = "11/11/2023 16:56"
$ws.Cells.Item($newRow, 14).Value = 4.43
$ws.Cells.Item($newRow, 15).Value = "05/11/2023 17:12"
$ws.Cells.Item($newRow, 16).Value = 5.13
$ws.Cells.Item($newRow, 17).Value = "11/11/2023 16:56"
$ws.Cells.Item($newRow, 18).Value = 1.46
$ws.Cells.Item($newRow, 19).Value = "05/11/2023 17:12"
$ws.Cells.Item($newRow, 20).Value = 1.26
$ws.Cells.Item($newRow, 21).Value = "11/11/2023 16:56"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/turkey/1-lig/altay-eyupspor/KM6278fA/"
